# "infection spreads within the regions"
#
# 1. Product Backlog sheet: add a new backlog item "Insects spread between
#    regions " with an estimate of 2, and move the total SUM formula down
#    one row (it now sums B2:B9 instead of B2:B8).
# 2. Sprint 2 sheet: extend the "remaining" burndown series (row 3) with two
#    more days of data (M3=4, N3=1).
# 3. Make "Product Backlog" the active sheet/tab again (it was "Sprint 2"),
#    updating each sheet's remembered selection along the way.

$wb = $excel.ActiveWorkbook

$wsBacklog = $wb.Worksheets.Item("Product Backlog")
$wsSprint2 = $wb.Worksheets.Item("Sprint 2")

# --- Product Backlog: restructure the total row and add the new item -----
# Remove the old total formula that lived on row 10.
$wsBacklog.Range("B10").ClearContents()

# New backlog entry on row 9.
$wsBacklog.Range("A9").Value = "Insects spread between regions "
$wsBacklog.Range("B9").Value = 2

# Recreate the running total one row down, now covering the new row.
$wsBacklog.Range("B11").Formula = "= SUM(B2:B9)"

# --- Sprint 2: two more days of burndown data on the "remaining" row -----
$wsSprint2.Range("M3").Value = 4
$wsSprint2.Range("N3").Value = 1

# --- Selections / active tab ---------------------------------------------
# Sprint 2 is no longer the selected tab; its remembered selection moves to B6.
[void]$wsSprint2.Activate()
$null = $wsSprint2.Range("B6").Select()

# Product Backlog becomes the selected tab again; remembered selection B12.
[void]$wsBacklog.Activate()
$null = $wsBacklog.Range("B12").Select()
